$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 24 : new activity log entry ---
$ws.Cells.Item(24, 2).Value = 779
$ws.Cells.Item(24, 3).Formula = '="5-4-2020"'
$ws.Cells.Item(24, 3).Copy()
$ws.Cells.Item(24, 3).PasteSpecial(-4163)
$ws.Cells.Item(24, 4).Value = 0.98402777777777783
$ws.Cells.Item(24, 4).NumberFormat = "h:mm AM/PM"
$ws.Cells.Item(24, 5).Value = 0.99097222222222225
$ws.Cells.Item(24, 7).Value = "Reviewed report together with team member for possible issues."

# --- Row 25 : new activity log entry ---
$ws.Cells.Item(25, 2).Value = 779
$ws.Cells.Item(25, 3).Formula = '="5-4-2020"'
$ws.Cells.Item(25, 3).Copy()
$ws.Cells.Item(25, 3).PasteSpecial(-4163)
$ws.Cells.Item(25, 4).Value = 0.99097222222222225
$ws.Cells.Item(25, 4).NumberFormat = "h:mm AM/PM"
$ws.Cells.Item(25, 5).Value = 0.99652777777777779
$ws.Cells.Item(25, 7).Value = "Revised report together with team member for clarity"

# --- Clear clipboard state from the copy operations above ---
$excel.CutCopyMode = 0

# --- Update the view: zoom to 70% and move the selection ---
$excel.ActiveWindow.Zoom = 70
$ws.Range("D29").Select()
